$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "September 21, 2025", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Split the mailing-address paragraph "2960 Sanor Pl, Santa Clara CA 95051"
#    (the copy in the letterhead, not the one repeated inside the info table)
#    into two separate paragraphs: "2960 Sanor Pl" and "Santa Clara, CA 95051".
#    Find/Replace only touches the first match when the replacement text
#    contains a new paragraph mark, which conveniently leaves the table's
#    copy of the address untouched.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2960 Sanor Pl, Santa Clara CA 95051", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2960 Sanor Pl^pSanta Clara, CA 95051", 2) | Out-Null

# The run created by the paragraph split above does not automatically pick up
# the Arial/11pt formatting used throughout the letter, so reapply it
# explicitly to the new "Santa Clara, CA 95051" paragraph. Locate it by
# indexed access (rather than re-enumerating $d.Paragraphs with foreach)
# since the paragraph count changed above.
$addrIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq "Santa Clara, CA 95051") {
        $addrIndex = $i
    }
}
if ($addrIndex -gt 0) {
    $addrRange = $d.Paragraphs.Item($addrIndex).Range
    $addrRange.Font.Name = "Arial"
    $addrRange.Font.NameAscii = "Arial"
    $addrRange.Font.NameOther = "Arial"
    $addrRange.Font.NameBi = "Arial"
    $addrRange.Font.Size = 11
    $addrRange.Font.SizeBi = 11
}

# ---------------------------------------------------------------------------
# 3. Remove the empty "No Spacing" paragraph that immediately follows the
#    "... Board of Directors" signature line.
# ---------------------------------------------------------------------------
$boardIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "Board of Directors") {
        $boardIndex = $i
    }
}
if ($boardIndex -gt 0) {
    $afterBoard = $d.Paragraphs.Item($boardIndex + 1)
    if ($afterBoard.Range.Text.TrimEnd([char]13, [char]7) -eq "") {
        $afterBoard.Range.Delete() | Out-Null
    }
}

Write-Host "Edit complete."
